$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 57.5
$ws.Range("I5").Value = 43.333332
$ws.Range("K5").Value = 43.333332
$ws.Range("M5").Value = 71.666668
# Row 18
$ws.Range("H18").Value = 1703.4
$ws.Range("I18").Value = 1270.4445
$ws.Range("K18").Value = 1270.4445
$ws.Range("M18").Value = -986.4445000000001
# Row 52
$ws.Range("H52").Value = 1491.8
$ws.Range("J52").Value = 9999
$ws.Range("L52").Value = 29997
$ws.Range("N52").Value = -30317
# Row 76
$ws.Range("H76").Value = 5895
$ws.Range("I76").Value = 4290
$ws.Range("K76").Value = 4290
$ws.Range("M76").Value = -3975
# Row 79
$ws.Range("H79").Value = 5895
$ws.Range("I79").Value = 4290
$ws.Range("K79").Value = 4290
$ws.Range("M79").Value = -3198
# Row 88
$ws.Range("H88").Value = 578327.8
$ws.Range("I88").Value = 1700
$ws.Range("J88").Value = 735589.9399999999
$ws.Range("K88").Value = 1700
$ws.Range("L88").Value = 735589.9399999999
$ws.Range("M88").Value = -1294
$ws.Range("N88").Value = -736401.9399999999
# Row 91
$ws.Range("H91").Value = 578327.8
$ws.Range("I91").Value = 1700
$ws.Range("J91").Value = 735589.9399999999
$ws.Range("K91").Value = 1700
$ws.Range("L91").Value = 735589.9399999999
$ws.Range("M91").Value = -296
$ws.Range("N91").Value = -738397.9399999999
# Row 92
$ws.Range("H92").Value = 333.55554
$ws.Range("I92").Value = 333.55554
$ws.Range("K92").Value = 333.55554
$ws.Range("M92").Value = 914.4444599999999
# Row 106
$ws.Range("H106").Value = 3611.0625
$ws.Range("I106").Value = 1753.7778
$ws.Range("K106").Value = 1753.7778
$ws.Range("M106").Value = -1122.7778
# Row 138
$ws.Range("H138").Value = 1877.61
$ws.Range("I138").Value = 1679.8
$ws.Range("J138").Value = 1899.5889
$ws.Range("K138").Value = 5039.4
$ws.Range("L138").Value = 5698.7667
$ws.Range("M138").Value = 100.6000000000004
$ws.Range("N138").Value = -15978.7667

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 583.6667
$ws.Range("I2").Value = 474.7619
$ws.Range("J2").Value = 964.8333
$ws.Range("K2").Value = 474.7619
$ws.Range("L2").Value = 964.8333
$ws.Range("M2").Value = -361.7619
$ws.Range("N2").Value = -1190.8333
# Row 32
$ws.Range("H32").Value = 10647571
$ws.Range("I32").Value = 13164638
$ws.Range("J32").Value = 19955.111
$ws.Range("K32").Value = 13164638
$ws.Range("L32").Value = 19955.111
$ws.Range("M32").Value = -13164351
$ws.Range("N32").Value = -20529.111
# Row 76
$ws.Range("H76").Value = 40287.332
$ws.Range("J76").Value = 40287.332
$ws.Range("L76").Value = 40287.332
$ws.Range("N76").Value = -40963.332
# Row 79
$ws.Range("H79").Value = 40287.332
$ws.Range("J79").Value = 40287.332
$ws.Range("L79").Value = 40287.332
$ws.Range("N79").Value = -42627.332
# Row 92
$ws.Range("H92").Value = 67747.25
$ws.Range("J92").Value = 67747.25
$ws.Range("L92").Value = 67747.25
$ws.Range("N92").Value = -72739.25
# Row 110
$ws.Range("H110").Value = 1953.4615
$ws.Range("I110").Value = 1807.4546
$ws.Range("K110").Value = 1807.4546
$ws.Range("M110").Value = 237.5454
# Row 116
$ws.Range("H116").Value = 583.6667
$ws.Range("I116").Value = 474.7619
$ws.Range("J116").Value = 964.8333
$ws.Range("K116").Value = 474.7619
$ws.Range("L116").Value = 964.8333
$ws.Range("M116").Value = 1819.2381
$ws.Range("N116").Value = -5552.8333
# Row 122
$ws.Range("H122").Value = 2240.6
$ws.Range("I122").Value = 1307.3572
$ws.Range("J122").Value = 3428.3635
$ws.Range("K122").Value = 3922.0716
$ws.Range("L122").Value = 10285.0905
$ws.Range("M122").Value = -1472.0716
$ws.Range("N122").Value = -15185.0905
# Row 125
$ws.Range("H125").Value = 73000
$ws.Range("J125").Value = 73000
$ws.Range("L125").Value = 73000
$ws.Range("N125").Value = -82840

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 583.6667
$ws.Range("I3").Value = 474.7619
$ws.Range("J3").Value = 964.8333
$ws.Range("K3").Value = 474.7619
$ws.Range("L3").Value = 964.8333
$ws.Range("M3").Value = -360.7619
$ws.Range("N3").Value = -1192.8333
# Row 134
$ws.Range("H134").Value = 911366.6
$ws.Range("I134").Value = 2446.5557
$ws.Range("J134").Value = 5001507
$ws.Range("K134").Value = 7339.6671
$ws.Range("L134").Value = 15004521
$ws.Range("M134").Value = -4804.6671
$ws.Range("N134").Value = -15009591

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 16787.625
$ws.Range("I6").Value = 859.6
$ws.Range("K6").Value = 859.6
$ws.Range("M6").Value = -746.6
# Row 74
$ws.Range("H74").Value = 48156.75
$ws.Range("J74").Value = 48156.75
$ws.Range("L74").Value = 48156.75
$ws.Range("N74").Value = -49904.75
# Row 77
$ws.Range("H77").Value = 48156.75
$ws.Range("J77").Value = 48156.75
$ws.Range("L77").Value = 144470.25
$ws.Range("N77").Value = -153206.25

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 786.1429000000001
$ws.Range("I7").Value = 367.33334
$ws.Range("J7").Value = 1100.25
$ws.Range("K7").Value = 1102.00002
$ws.Range("L7").Value = 3300.75
$ws.Range("M7").Value = -990.0000199999999
$ws.Range("N7").Value = -3524.75
# Row 74
$ws.Range("H74").Value = 13858.333
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 13858.333
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 41574.999
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -43696.999
# Row 77
$ws.Range("H77").Value = 13858.333
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 13858.333
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 124724.997
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -135332.997
# Row 115
$ws.Range("H115").Value = 15531.6
$ws.Range("I115").Value = 1209.3334
$ws.Range("K115").Value = 3628.0002
$ws.Range("M115").Value = -2453.0002
# Row 131
$ws.Range("H131").Value = 5799.7334
$ws.Range("J131").Value = 6032.884
$ws.Range("L131").Value = 18098.652
$ws.Range("N131").Value = -28178.652
# Row 137
$ws.Range("H137").Value = 7504.5
$ws.Range("I137").Value = 7010
$ws.Range("J137").Value = 7999
$ws.Range("K137").Value = 21030
$ws.Range("L137").Value = 23997
$ws.Range("M137").Value = -15930
$ws.Range("N137").Value = -34197
# Row 140
$ws.Range("H140").Value = 2432.9167
$ws.Range("I140").Value = 2408.261
$ws.Range("K140").Value = 7224.782999999999
$ws.Range("M140").Value = -2044.782999999999

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 1705.7222
$ws.Range("I122").Value = 1293.4375
$ws.Range("K122").Value = 3880.3125
$ws.Range("M122").Value = -1430.3125

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 45455228
$ws.Range("I55").Value = 52632236
$ws.Range("J55").Value = 856.6667
$ws.Range("K55").Value = 52632236
$ws.Range("L55").Value = 856.6667
$ws.Range("M55").Value = -52632063
$ws.Range("N55").Value = -1202.6667
# Row 122
$ws.Range("H122").Value = 5809.2856
$ws.Range("I122").Value = 5252.353
$ws.Range("K122").Value = 15757.059
$ws.Range("M122").Value = -13307.059

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 110249.75
$ws.Range("I2").Value = 110249.75
$ws.Range("K2").Value = 110249.75
$ws.Range("M2").Value = -110137.75
# Row 5
$ws.Range("H5").Value = 3586666.5
$ws.Range("I5").Value = 303999.8
$ws.Range("K5").Value = 303999.8
$ws.Range("M5").Value = -303887.8
# Row 18
$ws.Range("H18").Value = 8833.333000000001
$ws.Range("I18").Value = 8500
$ws.Range("K18").Value = 8500
$ws.Range("M18").Value = -8327
# Row 39
$ws.Range("H39").Value = 12277.223
$ws.Range("J39").Value = 30495
$ws.Range("L39").Value = 30495
$ws.Range("N39").Value = -31321
# Row 42
$ws.Range("H42").Value = 29999
$ws.Range("J42").Value = 29999
$ws.Range("L42").Value = 29999
$ws.Range("N42").Value = -30755
# Row 43
$ws.Range("H43").Value = 98013.5
$ws.Range("I43").Value = 96027
$ws.Range("K43").Value = 96027
$ws.Range("M43").Value = -95878
# Row 100
$ws.Range("H100").Value = 1224.2222
$ws.Range("I100").Value = 2017.6666
$ws.Range("J100").Value = 827.5
$ws.Range("K100").Value = 4035.3332
$ws.Range("L100").Value = 1655
$ws.Range("M100").Value = -3494.3332
$ws.Range("N100").Value = -2737

Write-Output "Applied Behemoth_Profits updates"